$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the MuSCs sending-cluster block (old rows 11-13)
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(11).Delete()

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Pspn"
$ws.Cells.Item(2, 3).Value = "Ret"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.5247883333333333
$ws.Cells.Item(2, 8).Value = 1.574365
$ws.Cells.Item(2, 9).Value = 0.1674845870648259
$ws.Cells.Item(2, 10).Value = 0.1674845870648259
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.219350333333333
$ws.Cells.Item(2, 14).Value = 3.658051
$ws.Cells.Item(2, 15).Value = 0.2422674834150417
$ws.Cells.Item(2, 16).Value = 0.2576244469655636
$ws.Cells.Item(2, 17).Value = 0.6399008291794444
$ws.Cells.Item(2, 18).Value = 5.759107462615
$ws.Cells.Item(2, 19).Value = 0.04057606941900281
$ws.Cells.Item(2, 20).Value = 0.04314812411783156

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Pspn"
$ws.Cells.Item(3, 3).Value = "Ret"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.5247883333333333
$ws.Cells.Item(3, 8).Value = 1.574365
$ws.Cells.Item(3, 9).Value = 0.1674845870648259
$ws.Cells.Item(3, 10).Value = 0.1674845870648259
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.913661333333334
$ws.Cells.Item(3, 14).Value = 8.740984000000001
$ws.Cells.Item(3, 15).Value = 0.5789028628226193
$ws.Cells.Item(3, 16).Value = 0.6155986258624717
$ws.Cells.Item(3, 17).Value = 1.529055475017778
$ws.Cells.Item(3, 18).Value = 13.76149927516
$ws.Cells.Item(3, 19).Value = 0.09695730693049191
$ws.Cells.Item(3, 20).Value = 0.1031032816502503

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Pspn"
$ws.Cells.Item(4, 3).Value = "Ret"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.5247883333333333
$ws.Cells.Item(4, 8).Value = 1.574365
$ws.Cells.Item(4, 9).Value = 0.1674845870648259
$ws.Cells.Item(4, 10).Value = 0.1674845870648259
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.900063
$ws.Cells.Item(4, 14).Value = 1.800126
$ws.Cells.Item(4, 15).Value = 0.178829653762339
$ws.Cells.Item(4, 16).Value = 0.1267769271719646
$ws.Cells.Item(4, 17).Value = 0.4723425616649999
$ws.Cells.Item(4, 18).Value = 2.83405536999
$ws.Cells.Item(4, 19).Value = 0.02995121071533113
$ws.Cells.Item(4, 20).Value = 0.02123318129674399

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Pspn"
$ws.Cells.Item(5, 3).Value = "Ret"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.239930333333333
$ws.Cells.Item(5, 8).Value = 6.719791
$ws.Cells.Item(5, 9).Value = 0.7148668960482055
$ws.Cells.Item(5, 10).Value = 0.7148668960482057
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.219350333333333
$ws.Cells.Item(5, 14).Value = 3.658051
$ws.Cells.Item(5, 15).Value = 0.2422674834150417
$ws.Cells.Item(5, 16).Value = 0.2576244469655636
$ws.Cells.Item(5, 17).Value = 2.731259798593444
$ws.Cells.Item(5, 18).Value = 24.581338187341
$ws.Cells.Item(5, 19).Value = 0.173189003882321
$ws.Cells.Item(5, 20).Value = 0.1841671887484081

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Pspn"
$ws.Cells.Item(6, 3).Value = "Ret"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.239930333333333
$ws.Cells.Item(6, 8).Value = 6.719791
$ws.Cells.Item(6, 9).Value = 0.7148668960482055
$ws.Cells.Item(6, 10).Value = 0.7148668960482057
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.913661333333334
$ws.Cells.Item(6, 14).Value = 8.740984000000001
$ws.Cells.Item(6, 15).Value = 0.5789028628226193
$ws.Cells.Item(6, 16).Value = 0.6155986258624717
$ws.Cells.Item(6, 17).Value = 6.526398401593778
$ws.Cells.Item(6, 18).Value = 58.73758561434401
$ws.Cells.Item(6, 19).Value = 0.413838492659426
$ws.Cells.Item(6, 20).Value = 0.4400710788818459

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Pspn"
$ws.Cells.Item(7, 3).Value = "Ret"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.239930333333333
$ws.Cells.Item(7, 8).Value = 6.719791
$ws.Cells.Item(7, 9).Value = 0.7148668960482055
$ws.Cells.Item(7, 10).Value = 0.7148668960482057
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.900063
$ws.Cells.Item(7, 14).Value = 1.800126
$ws.Cells.Item(7, 15).Value = 0.178829653762339
$ws.Cells.Item(7, 16).Value = 0.1267769271719646
$ws.Cells.Item(7, 17).Value = 2.016078415611
$ws.Cells.Item(7, 18).Value = 12.096470493666
$ws.Cells.Item(7, 19).Value = 0.1278393995064586
$ws.Cells.Item(7, 20).Value = 0.09062862841795175

# Row 8
$ws.Cells.Item(8, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value = "Pspn"
$ws.Cells.Item(8, 3).Value = "Ret"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.3686343333333333
$ws.Cells.Item(8, 8).Value = 1.105903
$ws.Cells.Item(8, 9).Value = 0.1176485168869685
$ws.Cells.Item(8, 10).Value = 0.1176485168869685
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.219350333333333
$ws.Cells.Item(8, 14).Value = 3.658051
$ws.Cells.Item(8, 15).Value = 0.2422674834150417
$ws.Cells.Item(8, 16).Value = 0.2576244469655636
$ws.Cells.Item(8, 17).Value = 0.4494943972281111
$ws.Cells.Item(8, 18).Value = 4.045449575053
$ws.Cells.Item(8, 19).Value = 0.02850241011371789
$ws.Cells.Item(8, 20).Value = 0.03030913409932403

# Row 9
$ws.Cells.Item(9, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value = "Pspn"
$ws.Cells.Item(9, 3).Value = "Ret"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.3686343333333333
$ws.Cells.Item(9, 8).Value = 1.105903
$ws.Cells.Item(9, 9).Value = 0.1176485168869685
$ws.Cells.Item(9, 10).Value = 0.1176485168869685
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.913661333333334
$ws.Cells.Item(9, 14).Value = 8.740984000000001
$ws.Cells.Item(9, 15).Value = 0.5789028628226193
$ws.Cells.Item(9, 16).Value = 0.6155986258624717
$ws.Cells.Item(9, 17).Value = 1.074075603172445
$ws.Cells.Item(9, 18).Value = 9.666680428552002
$ws.Cells.Item(9, 19).Value = 0.06810706323270133
$ws.Cells.Item(9, 20).Value = 0.0724242653303756

# Row 10
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Pspn"
$ws.Cells.Item(10, 3).Value = "Ret"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.3686343333333333
$ws.Cells.Item(10, 8).Value = 1.105903
$ws.Cells.Item(10, 9).Value = 0.1176485168869685
$ws.Cells.Item(10, 10).Value = 0.1176485168869685
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.900063
$ws.Cells.Item(10, 14).Value = 1.800126
$ws.Cells.Item(10, 15).Value = 0.178829653762339
$ws.Cells.Item(10, 16).Value = 0.1267769271719646
$ws.Cells.Item(10, 17).Value = 0.331794123963
$ws.Cells.Item(10, 18).Value = 1.990764743778
$ws.Cells.Item(10, 19).Value = 0.02103904354054927
$ws.Cells.Item(10, 20).Value = 0.01491511745726885
